# Scheduled runner update: refresh cached market-board price snapshots
# (currentAveragePrice*, Leve buy-cost, and derived profit columns) for
# the leves whose figures moved since the last sync, across all 8 sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H27").Value = 10000
$ws.Range("J27").Value = 10000
$ws.Range("L27").Value = 30000
$ws.Range("N27").Value = -30202

$ws.Range("H32").Value = 4600
$ws.Range("I32").Value = 4000
$ws.Range("J32").Value = 5000
$ws.Range("K32").Value = 4000
$ws.Range("L32").Value = 5000
$ws.Range("M32").Value = -3674
$ws.Range("N32").Value = -5652

$ws.Range("H62").Value = 5699.1665
$ws.Range("I62").Value = 4899.727
$ws.Range("J62").Value = 6955.4287
$ws.Range("K62").Value = 4899.727
$ws.Range("L62").Value = 6955.4287
$ws.Range("M62").Value = -4275.727
$ws.Range("N62").Value = -8203.4287

$ws.Range("H65").Value = 5699.1665
$ws.Range("I65").Value = 4899.727
$ws.Range("J65").Value = 6955.4287
$ws.Range("K65").Value = 24498.635
$ws.Range("L65").Value = 34777.14350000001
$ws.Range("M65").Value = -21378.635
$ws.Range("N65").Value = -41017.14350000001

$ws.Range("H92").Value = 38462644
$ws.Range("I92").Value = 41667760
$ws.Range("J92").Value = 1290
$ws.Range("K92").Value = 41667760
$ws.Range("L92").Value = 1290
$ws.Range("M92").Value = -41666512
$ws.Range("N92").Value = -3786

$ws.Range("H96").Value = 1489.5
$ws.Range("I96").Value = 1489.5
$ws.Range("K96").Value = 4468.5
$ws.Range("M96").Value = -3095.5

$ws.Range("H103").Value = 900.40625
$ws.Range("I103").Value = 537.7857
$ws.Range("K103").Value = 1613.3571
$ws.Range("M103").Value = -1027.3571

$ws.Range("H132").Value = 2588.1936
$ws.Range("I132").Value = 2015.0416
$ws.Range("J132").Value = 4553.2856
$ws.Range("K132").Value = 6045.1248
$ws.Range("L132").Value = 13659.8568
$ws.Range("M132").Value = -3515.1248
$ws.Range("N132").Value = -18719.8568

$ws.Range("H137").Value = 84981
$ws.Range("I137").Value = 84981
$ws.Range("K137").Value = 254943
$ws.Range("M137").Value = -252393

$ws.Range("H138").Value = 2636.5
$ws.Range("J138").Value = 2738.988
$ws.Range("L138").Value = 8216.964
$ws.Range("N138").Value = -18496.964

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 23967252
$ws.Range("I32").Value = 26582078
$ws.Range("K32").Value = 26582078
$ws.Range("M32").Value = -26581791

$ws.Range("H61").Value = 3757.3914
$ws.Range("I61").Value = 3152.4443
$ws.Range("J61").Value = 4146.2856
$ws.Range("K61").Value = 3152.4443
$ws.Range("L61").Value = 4146.2856
$ws.Range("M61").Value = -2940.4443
$ws.Range("N61").Value = -4570.2856

$ws.Range("H82").Value = 94995
$ws.Range("J82").Value = 94995
$ws.Range("L82").Value = 94995
$ws.Range("N82").Value = -95717

$ws.Range("H85").Value = 94995
$ws.Range("J85").Value = 94995
$ws.Range("L85").Value = 94995
$ws.Range("N85").Value = -97491

$ws.Range("H125").Value = 51201.75
$ws.Range("J125").Value = 51201.75
$ws.Range("L125").Value = 51201.75
$ws.Range("N125").Value = -61041.75

$ws.Range("H132").Value = 2576.0293
$ws.Range("I132").Value = 2372.9
$ws.Range("J132").Value = 4099.5
$ws.Range("K132").Value = 7118.700000000001
$ws.Range("L132").Value = 12298.5
$ws.Range("M132").Value = -4588.700000000001
$ws.Range("N132").Value = -17358.5

$ws.Range("H136").Value = 3757.3914
$ws.Range("I136").Value = 3152.4443
$ws.Range("J136").Value = 4146.2856
$ws.Range("K136").Value = 9457.332900000001
$ws.Range("L136").Value = 12438.8568
$ws.Range("M136").Value = -6907.332900000001
$ws.Range("N136").Value = -17538.8568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1076.55
$ws.Range("I94").Value = 1147.091
$ws.Range("K94").Value = 1147.091
$ws.Range("M94").Value = -696.0909999999999

$ws.Range("H105").Value = 3060
$ws.Range("I105").Value = 1855
$ws.Range("K105").Value = 1855
$ws.Range("M105").Value = -108

$ws.Range("H134").Value = 2749586.5
$ws.Range("I134").Value = 3107515.2
$ws.Range("K134").Value = 9322545.600000001
$ws.Range("M134").Value = -9320010.600000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 2449.5
$ws.Range("I41").Value = 2449.5
$ws.Range("K41").Value = 2449.5
$ws.Range("M41").Value = -2021.5

$ws.Range("H96").Value = 35833.332
$ws.Range("J96").Value = 35833.332
$ws.Range("L96").Value = 35833.332
$ws.Range("N96").Value = -41325.332

$ws.Range("H99").Value = 3249.5
$ws.Range("I99").Value = 2499.5
$ws.Range("J99").Value = 3499.5
$ws.Range("K99").Value = 2499.5
$ws.Range("L99").Value = 3499.5
$ws.Range("M99").Value = -1001.5
$ws.Range("N99").Value = -6495.5

$ws.Range("H126").Value = 3249.5
$ws.Range("I126").Value = 2499.5
$ws.Range("J126").Value = 3499.5
$ws.Range("K126").Value = 7498.5
$ws.Range("L126").Value = 10498.5
$ws.Range("M126").Value = -5028.5
$ws.Range("N126").Value = -15438.5

$ws.Range("H134").Value = 2432.077
$ws.Range("I134").Value = 1758.8889
$ws.Range("K134").Value = 5276.6667
$ws.Range("M134").Value = -2741.6667

$ws.Range("H141").Value = 490560.6
$ws.Range("J141").Value = 490560.6
$ws.Range("L141").Value = 490560.6
$ws.Range("N141").Value = -500920.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 4998.3335
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()

$ws.Range("H40").Value = 2373.9
$ws.Range("J40").Value = 5463.5
$ws.Range("L40").Value = 21854
$ws.Range("N40").Value = -21992

$ws.Range("H46").Value = 3449.6667
$ws.Range("J46").Value = 4999.5
$ws.Range("L46").Value = 14998.5
$ws.Range("N46").Value = -15180.5

$ws.Range("H107").Value = 877.5484
$ws.Range("J107").Value = 642.58826
$ws.Range("L107").Value = 1927.76478
$ws.Range("N107").Value = -5767.76478

$ws.Range("H121").Value = 4769062
$ws.Range("I121").Value = 955.25
$ws.Range("J121").Value = 6040557
$ws.Range("K121").Value = 2865.75
$ws.Range("L121").Value = 18121671
$ws.Range("M121").Value = -1555.75
$ws.Range("N121").Value = -18124291

$ws.Range("H133").Value = 3190.9412
$ws.Range("I133").Value = 2437.1667
$ws.Range("K133").Value = 7311.500100000001
$ws.Range("M133").Value = -2251.500100000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 1050.5
$ws.Range("I10").Value = 2001.5
$ws.Range("J10").Value = 99.5
$ws.Range("K10").Value = 2001.5
$ws.Range("L10").Value = 99.5
$ws.Range("M10").Value = -1832.5
$ws.Range("N10").Value = -437.5

$ws.Range("H62").Value = 82000
$ws.Range("I62").Value = 48000
$ws.Range("J62").Value = 116000
$ws.Range("K62").Value = 48000
$ws.Range("L62").Value = 116000
$ws.Range("M62").Value = -47314
$ws.Range("N62").Value = -117372

$ws.Range("H65").Value = 82000
$ws.Range("I65").Value = 48000
$ws.Range("J65").Value = 116000
$ws.Range("K65").Value = 144000
$ws.Range("L65").Value = 348000
$ws.Range("M65").Value = -140568
$ws.Range("N65").Value = -354864

$ws.Range("H70").Value = 19690.35
$ws.Range("J70").Value = 4488.8335
$ws.Range("L70").Value = 4488.8335
$ws.Range("N70").Value = -5028.8335

$ws.Range("H73").Value = 19690.35
$ws.Range("J73").Value = 4488.8335
$ws.Range("L73").Value = 4488.8335
$ws.Range("N73").Value = -6360.8335

$ws.Range("H80").Value = 2545.818
$ws.Range("I80").Value = 2410.4
$ws.Range("J80").Value = 3900
$ws.Range("K80").Value = 2410.4
$ws.Range("L80").Value = 3900
$ws.Range("M80").Value = -1412.4
$ws.Range("N80").Value = -5896

$ws.Range("H83").Value = 2545.818
$ws.Range("I83").Value = 2410.4
$ws.Range("J83").Value = 3900
$ws.Range("K83").Value = 12052
$ws.Range("L83").Value = 19500
$ws.Range("M83").Value = -7060
$ws.Range("N83").Value = -29484

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4735
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()

$ws.Range("H81").Value = 69990
$ws.Range("J81").Value = 69990
$ws.Range("L81").Value = 69990
$ws.Range("N81").Value = -71986

$ws.Range("H82").Value = 1483.3334
$ws.Range("I82").Value = 1370.5555
$ws.Range("J82").Value = 1652.5
$ws.Range("K82").Value = 1370.5555
$ws.Range("L82").Value = 1652.5
$ws.Range("M82").Value = -1009.5555
$ws.Range("N82").Value = -2374.5

$ws.Range("H84").Value = 69990
$ws.Range("J84").Value = 69990
$ws.Range("L84").Value = 209970
$ws.Range("N84").Value = -219954

$ws.Range("H85").Value = 1483.3334
$ws.Range("I85").Value = 1370.5555
$ws.Range("J85").Value = 1652.5
$ws.Range("K85").Value = 1370.5555
$ws.Range("L85").Value = 1652.5
$ws.Range("M85").Value = -122.5554999999999
$ws.Range("N85").Value = -4148.5

$ws.Range("H122").Value = 21994.4
$ws.Range("I122").Value = 21993.125
$ws.Range("J122").Value = 21999.5
$ws.Range("K122").Value = 65979.375
$ws.Range("L122").Value = 65998.5
$ws.Range("M122").Value = -63529.375
$ws.Range("N122").Value = -70898.5

$ws.Range("H132").Value = 4283.7407
$ws.Range("I132").Value = 4003.818
$ws.Range("J132").Value = 5515.4
$ws.Range("K132").Value = 12011.454
$ws.Range("L132").Value = 16546.2
$ws.Range("M132").Value = -9481.454000000002
$ws.Range("N132").Value = -21606.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2533.25
$ws.Range("I81").Value = 1551.2307
$ws.Range("K81").Value = 3102.4614
$ws.Range("M81").Value = -2041.4614

$ws.Range("H84").Value = 2533.25
$ws.Range("I84").Value = 1551.2307
$ws.Range("K84").Value = 15512.307
$ws.Range("M84").Value = -10208.307
